# Fruta / hortaliza, semanal
# Insert a new data row at row 134 (pushing existing rows 134-176 down to
# 135-177) and populate it with the new weekly price-report entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 134..176 down by one to make room for the new record.
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new record's data.
$ws.Cells.Item(134, 1).Value  = 7
$ws.Cells.Item(134, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(134, 3).Value  = "Ñuble"
$ws.Cells.Item(134, 4).Value  = 45229
$ws.Cells.Item(134, 5).Value  = 16
$ws.Cells.Item(134, 6).Value  = 100112031
$ws.Cells.Item(134, 7).Value  = "Poroto verde"
$ws.Cells.Item(134, 8).Value  = "Sin especificar"
$ws.Cells.Item(134, 9).Value  = "Primera"
$ws.Cells.Item(134, 10).Value = 100
$ws.Cells.Item(134, 11).Value = 20000
$ws.Cells.Item(134, 12).Value = 20000
$ws.Cells.Item(134, 13).Value = 20000
$ws.Cells.Item(134, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(134, 15).Value = "Perú"
$ws.Cells.Item(134, 16).Value = 800
$ws.Cells.Item(134, 17).Value = 25
$ws.Cells.Item(134, 18).Value = "Hortaliza"
